# Refresh the traffic-incident table with speedband/incident data retrieved
# from the LTA (Land Transport Authority) DataMall API for the original route.
#
# - Row 2 now reports a "Vehicle breakdown" (was "Accident") with updated
#   coordinates/node id.
# - Rows 3-25 are refreshed "Roadwork" entries with newer lat/long/node values.
# - A new "Roadwork" entry is appended as row 26, growing the sheet to A1:D26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Type, Latitude, Longitude, Node
$incidents = @(
    ,@(2, "Vehicle breakdown", 1.424653170043902, 103.7910440091601, 440588376)
    ,@(3, "Roadwork", 1.329566083085618, 103.8398779449462, 1223810685)
    ,@(4, "Roadwork", 1.294395580676698, 103.8665957683079, 26782016)
    ,@(5, "Roadwork", 1.383978872974393, 103.9187785015995, 5726683755)
    ,@(6, "Roadwork", 1.399693726009656, 103.8580241805984, 1781710554)
    ,@(7, "Roadwork", 1.329431928728425, 103.8583187295057, 250830816)
    ,@(8, "Roadwork", 1.390077478097602, 103.7465491295866, 1734061685)
    ,@(9, "Roadwork", 1.388890087119251, 103.7578528352893, 5045819346)
    ,@(10, "Roadwork", 1.4008676440462, 103.8597329453334, 1781710559)
    ,@(11, "Roadwork", 1.387215261909439, 103.7562256569732, 236726630)
    ,@(12, "Roadwork", 1.330328271085954, 103.9152015923327, 1849458129)
    ,@(13, "Roadwork", 1.309615816817113, 103.8757620402856, 3066261840)
    ,@(14, "Roadwork", 1.300410885561765, 103.8779458712998, 1637064322)
    ,@(15, "Roadwork", 1.395616215146429, 103.8576328064338, 5158040655)
    ,@(16, "Roadwork", 1.298178361107145, 103.8003640080351, 247331373)
    ,@(17, "Roadwork", 1.333070893064765, 103.8890100726651, 239833041)
    ,@(18, "Roadwork", 1.335214775723557, 103.861906736975, 137302917)
    ,@(19, "Roadwork", 1.325673028980316, 103.7406708548471, 1350302489)
    ,@(20, "Roadwork", 1.337821500296478, 103.9805201027967, 4748287569)
    ,@(21, "Roadwork", 1.324029203005592, 103.7433344949127, 566531128)
    ,@(22, "Roadwork", 1.426491921158762, 103.7821463731425, 1110672475)
    ,@(23, "Roadwork", 1.304490794689495, 103.770414048807, 247184656)
    ,@(24, "Roadwork", 1.408640772693746, 103.8080272748565, 148992674)
    ,@(25, "Roadwork", 1.394240890722329, 103.7742911599704, 158014654)
    ,@(26, "Roadwork", 1.352252843616421, 103.6957155183491, 567939126)
)

foreach ($incident in $incidents) {
    $row = $incident[0]
    $ws.Cells.Item($row, 1).Value = $incident[1]   # Type
    $ws.Cells.Item($row, 2).Value = $incident[2]   # Latitude
    $ws.Cells.Item($row, 3).Value = $incident[3]   # Longitude
    $ws.Cells.Item($row, 4).Value = $incident[4]   # Node
}
